$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the Halkbank (column H) cells that were removed from this benchmark run.
$ws.Range("H3").Value = ""
$ws.Range("H4").Value = ""
$ws.Range("H5").Value = ""
$ws.Range("H6").Value = ""
$ws.Range("H8").Value = ""
$ws.Range("H9").Value = ""
$ws.Range("H10").Value = ""
$ws.Range("H11").Value = ""
$ws.Range("H13").Value = ""
$ws.Range("H14").Value = ""

# Updated benchmark figures for row 13 (GELEN SWIFT).
$ws.Range("D13").Value = "Hesaba: Asgari 1 TL | Azami 300 TL"
$ws.Range("K13").Value = "Hesaba: Asgari 1 TL | Azami 53,19 TL"
